$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$nl = [char]10

# Developer name
$ws.Range("C3").Value = "Hudson Drozdowski "

# Reusable Method Inputs blocks
$inputsValid      = 'color = "Blue"' + $nl + 'length = 5' + $nl + 'width = 5'
$inputsBlankColor = 'color = "     "' + $nl + 'length = 5' + $nl + 'width = 5'
$inputsBadLength  = 'color = "Blue"' + $nl + 'length = "Invalid Input"' + $nl + 'width = 5'
$inputsBadWidth   = 'color = "Blue"' + $nl + 'length = 5' + $nl + 'width = "Invalid Input"'

$none = "None"
$valueError = "ValueError()"

# Row 7 - __init__ / Attribute set to input values.
$ws.Range("E7").Value = $inputsValid
$ws.Range("F7").Value = $none
$ws.Range("G7").Value = "Object is initialized correctly."

# Row 8 - __init__ / Exception raised when color is blank
$ws.Range("E8").Value = $inputsBlankColor
$ws.Range("F8").Value = $none
$ws.Range("G8").Value = $valueError

# Row 9 - __init__ / Exception raised when length is not an integer.
$ws.Range("E9").Value = $inputsBadLength
$ws.Range("F9").Value = $none
$ws.Range("G9").Value = $valueError

# Row 10 - __init__ / Exception raised when width is not an integer.
$ws.Range("E10").Value = $inputsBadWidth
$ws.Range("F10").Value = $none
$ws.Range("G10").Value = $valueError

# Row 11 - __str__ / Returns string formatted appropriately
$ws.Range("E11").Value = $inputsValid
$ws.Range("F11").Value = $none
$ws.Range("G11").Value = '"The shape color is Blue. ' + $nl + 'This rectangle has four sides with the lengths of 5, 5, 5 and 5 centimeters."'

# Row 12 - calculate_area / Returns correct calculated value.
$ws.Range("E12").Value = $inputsValid
$ws.Range("F12").Value = $none
$ws.Range("G12").Value = 25

# Row 13 - calculate_perimeter / Returns correct calculated value.
$ws.Range("E13").Value = $inputsValid
$ws.Range("F13").Value = $none
$ws.Range("G13").Value = 20
